# Applies the edits described by the commit:
# "Added generic methods to read common data properties, read and write
#  data into excel file"
#
# Concrete data changes in the workbook:
#   - createUser!D2   : "sagardambal" -> "Ocean"
#   - InvalidData!D11 : "8%"          -> "Samsung galaxy s25 ultra"
#   - validateCred!A1 : "Admin"       -> "Modi"
#   - Active sheet / selection moves from validateCred (no selection) to
#     createUser, with the active cell at H9.

$wb = $excel.ActiveWorkbook

$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsInvalidData = $wb.Worksheets.Item("InvalidData")
$wsValidateCred = $wb.Worksheets.Item("validateCred")

# createUser: password value for row 2 updated
$wsCreateUser.Range("D2").Value = "Ocean"

# InvalidData: row 11 discount cell overwritten with new text, clearing the
# percentage number format it used to carry
$wsInvalidData.Range("D11").ClearFormats()
$wsInvalidData.Range("D11").Value = "Samsung galaxy s25 ultra"
$wsInvalidData.Columns.Item(4).AutoFit()

# validateCred: first data row's username changed
$wsValidateCred.Range("A1").Value = "Modi"

# Active sheet becomes createUser again, with H9 selected
$wsCreateUser.Activate()
$wsCreateUser.Range("H9").Select()
